# Applies the Assignment-2 roster table updates:
#  - merge the split "C" / "0891502" runs for Olumayowa Mosury's Student ID
#    into a single "C0891502" run
#  - fill in the blank roster row with Rafaeatul Kabir / C0888535 / 6
#
# Each target paragraph is replaced wholesale (via Range.InsertXML) with a
# fully-formed <w:p> that reproduces the original paragraph's identity
# attributes (paraId/textId/rsid*) and run formatting (italic, 28 half-pt /
# 14pt sizes for both Western and complex-script fonts), so the only visible
# change is the run/text content.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- 1) Olumayowa Mosury's Student ID: merge "C" + "0891502" -> "C0891502" ---
$studentIdCell = $tbl.Cell(3, 2)
$studentIdRange = $studentIdCell.Range
$studentIdRange.Collapse(1)
$studentIdXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="470FC415" w14:textId="4023208C" w:rsidR="003A63C3" w:rsidRDefault="001932C2" w:rsidP="008A490D">' +
  '<w:pPr><w:spacing w:after="480"/><w:jc w:val="center"/><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="001932C2"><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>C0891502</w:t></w:r>' +
  '</w:p>'
$studentIdRange.InsertXML($studentIdXml)

# --- 2) Blank roster row: Name = "Rafaeatul Kabir" ---
$nameCell = $tbl.Cell(4, 1)
$nameRange = $nameCell.Range
$nameRange.Collapse(1)
$nameXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="21D659BC" w14:textId="77777777" w:rsidR="009A159C" w:rsidRDefault="009A159C" w:rsidP="008A490D">' +
  '<w:pPr><w:spacing w:after="480"/><w:jc w:val="center"/><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Rafaeatul</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Kabir</w:t></w:r>' +
  '</w:p>'
$nameRange.InsertXML($nameXml)

# --- 3) Blank roster row: Student ID = "C" + "0888535" ---
$newStudentIdCell = $tbl.Cell(4, 2)
$newStudentIdRange = $newStudentIdCell.Range
$newStudentIdRange.Collapse(1)
$newStudentIdXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="373D4A0E" w14:textId="77777777" w:rsidR="009A159C" w:rsidRDefault="009A159C" w:rsidP="008A490D">' +
  '<w:pPr><w:spacing w:after="480"/><w:jc w:val="center"/><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>C</w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>0888535</w:t></w:r>' +
  '</w:p>'
$newStudentIdRange.InsertXML($newStudentIdXml)

# --- 4) Blank roster row: Group# = "6" ---
$groupCell = $tbl.Cell(4, 3)
$groupRange = $groupCell.Range
$groupRange.Collapse(1)
$groupXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="454EBCDF" w14:textId="77777777" w:rsidR="009A159C" w:rsidRDefault="009A159C" w:rsidP="008A490D">' +
  '<w:pPr><w:spacing w:after="480"/><w:jc w:val="center"/><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>6</w:t></w:r>' +
  '</w:p>'
$groupRange.InsertXML($groupXml)
